$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '61.701.71'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.83%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.897.71'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.07%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '571.68'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.84%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.68'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.45%  '
$ws.Range('E7').Value = '  +0.02%  '
$ws.Range('E8').Value = '  -0.91%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.897.74'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.06%  '
$ws.Range('E10').Value = '  -7.43%  '
$ws.Range('E11').Value = '  -1.90%  '
$ws.Range('E12').Value = '  -2.74%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000235'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -2.09%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '32.06'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.21%  '
$ws.Range('E15').Value = '  -0.81%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.384.42'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.04%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.733.07'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.74%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.62'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.59%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.904.94'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.13%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '435.05'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -1.64%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.24'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.26%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.656'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.08%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.89'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.93%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '79.16'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.97%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.93'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.75%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.18'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -9.98%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.04'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000107'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +10.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.10'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.52'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.06%  '
$ws.Range('E32').Value = '  -3.88%  '
$ws.Range('B33').Value = 'FirstDigitalUSD'
$ws.Range('C33').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.28%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.107'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.70%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '25.61'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.49%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.960'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -3.60%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.44'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.85%  '
$ws.Range('B38').Value = 'OKB'
$ws.Range('C38').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '49.04'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.98%  '
$ws.Range('B39').Value = 'dogwifhat'
$ws.Range('C39').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.91'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.96%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.95'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.64%  '
$ws.Range('E41').Value = '  -1.08%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.26'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.269'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -4.43%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '38.70'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -5.28%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.689.89'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.80%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '133.07'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.21%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0333'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -1.82%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '337.79'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -7.82%  '
$ws.Range('E50').Value = '  -2.06%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '21.63'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.79%  '
